$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The geolocation strings in column BL (rows 2-157) were stored wrapped in
# parentheses, e.g. "(36.1576629645, -86.8057331508)". Strip the parentheses
# so the cell holds the bare "lat, long" text (start of splitting these into
# usable coordinates for zip/tract lookups).
$values = @(
    "36.1576629645, -86.8057331508",
    "36.2758406196, -86.7227516725",
    "36.0133154771, -86.7176062426",
    "36.2962421445, -86.884373603",
    "36.1231320505, -86.7142265091",
    "36.1860223442, -86.7643746828",
    "36.0304702133, -86.7544989419",
    "36.2212110254, -86.7946381366",
    "36.1342444257, -86.9660257669",
    "36.2104122971, -86.6882030885",
    "36.0851071347, -86.6579014347",
    "36.1454955363, -86.7638145949",
    "36.120319646, -86.7820995582",
    "36.1286948607, -86.7851320823",
    "36.2021376395, -86.7609246036",
    "36.2156183387, -86.922992251",
    "36.1215972147, -86.8291852595",
    "36.0972128283, -86.5907016995",
    "36.1006385723, -86.7268837685",
    "36.2947153831, -86.6911187036",
    "36.0886699253, -86.6809974727",
    "36.0702023848, -86.6951799025",
    "36.042374025, -86.7688478485",
    "36.0627541445, -86.6647320159",
    "36.1457801977, -86.8050873817",
    "36.0931627912, -86.8736555705",
    "36.2229163853, -86.7598911314",
    "36.1711531762, -86.752019091",
    "36.2059653483, -86.7961197587",
    "36.1709368867, -86.644421803",
    "36.0593607088, -86.7364092465",
    "36.1714412594, -86.8285203739",
    "36.2488660642, -86.8538782901",
    "36.1911511856, -86.8085327041",
    "36.1790862802, -86.7590618024",
    "36.0521853395, -86.9407121165",
    "35.9996029242, -86.6586366088",
    "36.1804768223, -86.7383080372",
    "36.1106881287, -86.8763668451",
    "36.0657786763, -86.7267065082",
    "36.2509916244, -86.731871255",
    "36.1851740179, -86.6053672734",
    "36.0995287373, -86.7996847381",
    "36.1612468521, -86.8156490258",
    "36.0532072157, -86.6572040786",
    "36.344971962, -86.7677653742",
    "36.0587984091, -86.6901293191",
    "36.1439613559, -86.8662660609",
    "36.2230900506, -86.6632138216",
    "36.0848148306, -86.7119860954",
    "36.0810876055, -86.6188609306",
    "36.1133065622, -86.7274177563",
    "36.1068258846, -86.7023156259",
    "36.1728286331, -86.6903581821",
    "36.2111454796, -86.7293550121",
    "36.1836001359, -86.8846679987",
    "36.0899457741, -86.7174801595",
    "36.1125202586, -86.7940986863",
    "36.0609096698, -86.8696928139",
    "36.2287255904, -86.6339493996",
    "36.1360963098, -86.7411039245",
    "36.1446955813, -86.8426724912",
    "36.175652119, -86.814804507",
    "36.1597837914, -86.8530700515",
    "36.0383226992, -86.6692426787",
    "36.2343195911, -86.7732117462",
    "36.0679630565, -86.8949107152",
    "36.0506682578, -86.6371112264",
    "36.1523856593, -86.6522822722",
    "36.1457359079, -86.6254409227",
    "36.1460702635, -86.8156944707",
    "36.1719689182, -86.6206020339",
    "36.1399427822, -86.7819906718",
    "36.0971117892, -86.7878147099",
    "36.1701688187, -86.7363924126",
    "36.0873432616, -86.932612223",
    "36.194381526, -86.6421905953",
    "36.0192103596, -86.6464353266",
    "36.1329172709, -86.8300431737",
    "36.1839455221, -86.5833021676",
    "36.1195178608, -86.7427999612",
    "36.1346828403, -86.7023819538",
    "36.1834698574, -86.7113792554",
    "36.0753128653, -86.7434576848",
    "36.0517956287, -86.7048611505",
    "36.253365324, -86.6817557329",
    "36.0460647087, -86.9719220455",
    "36.2652149217, -86.7586096067",
    "36.1562850243, -86.6793913918",
    "36.0707379004, -86.8228348529",
    "36.2882514084, -86.7728204615",
    "36.1777432202, -86.790108485",
    "36.1017022136, -86.6046378648",
    "36.1276431581, -86.8678623837",
    "36.2246054475, -86.7196053106",
    "36.1503310917, -86.8286410961",
    "36.069945042, -86.7140695527",
    "36.1757150934, -86.8056204212",
    "36.0513249553, -86.724258158",
    "36.15764352, -86.7816792716",
    "36.0300312489, -86.6145110967",
    "36.1267253014, -86.634377235",
    "36.0715952616, -86.9578354061",
    "36.204125401, -86.7098381564",
    "36.190570195, -86.855339398",
    "36.1082517917, -86.8268057504",
    "36.197632419, -86.6144973804",
    "36.148106453, -86.5856887267",
    "36.2729362225, -86.6610202749",
    "36.1153457153, -86.8913814952",
    "36.1431909387, -86.7899130999",
    "36.1350470326, -86.8100174335",
    "36.1582500552, -86.8715276437",
    "36.073922218, -86.6732926633",
    "36.1363692187, -86.7683698189",
    "36.2269478533, -86.6077806723",
    "36.049231748, -86.6797916232",
    "36.2557501251, -86.7050849585",
    "36.3479505187, -86.8675409005",
    "36.0401325074, -87.0176548329",
    "36.0842376102, -86.7415232924",
    "36.2493399768, -86.7638027805",
    "36.0696879457, -86.6150485034",
    "36.124158094, -86.7550707151",
    "36.1895940492, -86.7439426762",
    "36.1353822457, -86.7150553656",
    "36.151096554, -86.7597078621",
    "36.1871201603, -86.7719899728",
    "36.0320124181, -86.7277054628",
    "36.0655938524, -86.9777719034",
    "36.0524304194, -86.8041894006",
    "36.0999620431, -86.639123776",
    "36.0831866378, -86.5615305811",
    "36.2059453952, -86.5921345588",
    "36.1012139547, -86.7435789922",
    "36.1247951937, -86.8014921222",
    "36.0712760339, -86.9382627391",
    "36.2249944111, -86.7399718241",
    "36.1711444309, -86.768865228",
    "36.1538311711, -86.7221674827",
    "36.1514089398, -86.8836263472",
    "36.166932328, -86.8036510368",
    "36.0960612808, -86.8261504838",
    "36.2681195238, -86.6959751266",
    "36.1639335636, -86.8274060699",
    "36.0722889795, -86.6476420172",
    "36.0973117939, -86.908717199",
    "36.1401682761, -86.7966089777",
    "36.2537083567, -86.803607107",
    "36.1985224304, -86.7337794886",
    "36.0449466794, -86.7359377605",
    "36.2080299864, -86.8352588264",
    "36.0609917887, -86.7571314086",
    "36.0769982456, -86.7054436498",
    "36.0360852644, -86.6924294734",
    "36.2432209865, -86.7170424941"
)

$startRow = 2
$col = 64  # column BL
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($startRow + $i, $col).Value = $values[$i]
}

# Update the sheet view: scroll/select so column BL (the last, newly-edited
# column) is in focus, matching the author switching their cursor there.
$ws.Activate()
$ws.Range("BL1:BL1048576").Select()
